$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "image" column
$ws.Range("E1").Value = "image"

# Image URLs used across the attendance rows
$img1 = "https://img.freepik.com/free-photo/handsome-young-man-with-arms-crossed-white-background_23-2148222620.jpg"
$img2 = "https://img.freepik.com/free-photo/handsome-confident-smiling-man-with-hands-crossed-chest_176420-18743.jpg"
$img3 = "https://img.freepik.com/free-photo/attractive-mixed-race-male-with-positive-smile-shows-white-teeth-keeps-hands-stomach-being-high-spirit-wears-white-shirt-rejoices-positive-moments-life-people-emotions-concept_273609-15527.jpg"
$img4 = "https://img.freepik.com/free-photo/confident-handsome-guy-posing-against-white-wall_176420-32936.jpg"
$img5 = "https://img.freepik.com/free-photo/fashionable-stylish-man-with-dark-eyes-casual-clothes-looking-aside-with-placid-thoughtful-look-pensive-guy-with-puzzled-expression-thinking-about-something-building-plans_176420-10331.jpg"
$img6 = "https://img.freepik.com/free-photo/thoughtful-concerned-man-thinking-trying-find-solution_176420-19574.jpg"
$img7 = "https://img.freepik.com/free-photo/serious-thoughtful-man-squinting-skeptical-thinking-as-making-choice_176420-19020.jpg"

# Fill column E top-to-bottom so new shared-string entries land in row order.
# Rows 2, 4 and 7 get a real clickable hyperlink (styled with the Hyperlink
# cell style); the rest just carry the plain image URL text.
$ws.Hyperlinks.Add($ws.Range("E2"), $img1)
$ws.Range("E3").Value = $img2
$ws.Hyperlinks.Add($ws.Range("E4"), $img3)
$ws.Range("E5").Value = $img4
$ws.Range("E6").Value = $img5
$ws.Hyperlinks.Add($ws.Range("E7"), $img6)
$ws.Range("E8").Value = $img7
$ws.Range("E9").Value = $img3
$ws.Range("E10").Value = $img6
$ws.Range("E11").Value = $img1

# Adjust the selection to mirror the post-edit active range
$ws.Range("E1:E11").Select()
